$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt7b"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.072919
$ws.Range("H2").Value = 0.218757
$ws.Range("I2").Value = 0.1477750351608889
$ws.Range("J2").Value = 0.1477750351608889
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 29.75868033333333
$ws.Range("N2").Value = 89.27604099999999
$ws.Range("O2").Value = 0.4948552779010537
$ws.Range("P2").Value = 0.4948552779010535
$ws.Range("Q2").Value = 2.169973211226333
$ws.Range("R2").Value = 19.529758901037
$ws.Range("S2").Value = 0.07312725609137966
$ws.Range("T2").Value = 0.07312725609137964

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt7b"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.072919
$ws.Range("H3").Value = 0.218757
$ws.Range("I3").Value = 0.1477750351608889
$ws.Range("J3").Value = 0.1477750351608889
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.55525033333333
$ws.Range("N3").Value = 52.665751
$ws.Range("O3").Value = 0.2919251856942525
$ws.Range("P3").Value = 0.2919251856942524
$ws.Range("Q3").Value = 1.280111299056333
$ws.Range("R3").Value = 11.521001691507
$ws.Range("S3").Value = 0.04313925458031719
$ws.Range("T3").Value = 0.04313925458031718

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt7b"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.072919
$ws.Range("H4").Value = 0.218757
$ws.Range("I4").Value = 0.1477750351608889
$ws.Range("J4").Value = 0.1477750351608889
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1421396666666667
$ws.Range("N4").Value = 0.426419
$ws.Range("O4").Value = 0.002363631836533717
$ws.Range("P4").Value = 0.002363631836533717
$ws.Range("Q4").Value = 0.01036468235366667
$ws.Range("R4").Value = 0.093282141183
$ws.Range("S4").Value = 0.0003492857777511665
$ws.Range("T4").Value = 0.0003492857777511664

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt7b"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.072919
$ws.Range("H5").Value = 0.218757
$ws.Range("I5").Value = 0.1477750351608889
$ws.Range("J5").Value = 0.1477750351608889
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.42872866666667
$ws.Range("N5").Value = 37.286186
$ws.Range("O5").Value = 0.2066765699758167
$ws.Range("P5").Value = 0.2066765699758166
$ws.Range("Q5").Value = 0.9062904656446666
$ws.Range("R5").Value = 8.156614190802001
$ws.Range("S5").Value = 0.03054163739510822
$ws.Range("T5").Value = 0.03054163739510822

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt7b"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.072919
$ws.Range("H6").Value = 0.218757
$ws.Range("I6").Value = 0.1477750351608889
$ws.Range("J6").Value = 0.1477750351608889
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.251329
$ws.Range("N6").Value = 0.753987
$ws.Range("O6").Value = 0.004179334592343558
$ws.Range("P6").Value = 0.004179334592343557
$ws.Range("Q6").Value = 0.018326659351
$ws.Range("R6").Value = 0.164939934159
$ws.Range("S6").Value = 0.0006176013163326885
$ws.Range("T6").Value = 0.0006176013163326884

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt7b"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.420527
$ws.Range("H7").Value = 1.261581
$ws.Range("I7").Value = 0.852224964839111
$ws.Range("J7").Value = 0.852224964839111
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 29.75868033333333
$ws.Range("N7").Value = 89.27604099999999
$ws.Range("O7").Value = 0.4948552779010537
$ws.Range("P7").Value = 0.4948552779010535
$ws.Range("Q7").Value = 12.51432856453567
$ws.Range("R7").Value = 112.628957080821
$ws.Range("S7").Value = 0.4217280218096739
$ws.Range("T7").Value = 0.4217280218096739

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt7b"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.420527
$ws.Range("H8").Value = 1.261581
$ws.Range("I8").Value = 0.852224964839111
$ws.Range("J8").Value = 0.852224964839111
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.55525033333333
$ws.Range("N8").Value = 52.665751
$ws.Range("O8").Value = 0.2919251856942525
$ws.Range("P8").Value = 0.2919251856942524
$ws.Range("Q8").Value = 7.382456756925667
$ws.Range("R8").Value = 66.44211081233101
$ws.Range("S8").Value = 0.2487859311139353
$ws.Range("T8").Value = 0.2487859311139352

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt7b"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.420527
$ws.Range("H9").Value = 1.261581
$ws.Range("I9").Value = 0.852224964839111
$ws.Range("J9").Value = 0.852224964839111
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1421396666666667
$ws.Range("N9").Value = 0.426419
$ws.Range("O9").Value = 0.002363631836533717
$ws.Range("P9").Value = 0.002363631836533717
$ws.Range("Q9").Value = 0.05977356760433334
$ws.Range("R9").Value = 0.537962108439
$ws.Range("S9").Value = 0.002014346058782551
$ws.Range("T9").Value = 0.00201434605878255

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt7b"
$ws.Range("C10").Value = "Fzd4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.420527
$ws.Range("H10").Value = 1.261581
$ws.Range("I10").Value = 0.852224964839111
$ws.Range("J10").Value = 0.852224964839111
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.42872866666667
$ws.Range("N10").Value = 37.286186
$ws.Range("O10").Value = 0.2066765699758167
$ws.Range("P10").Value = 0.2066765699758166
$ws.Range("Q10").Value = 5.226615980007334
$ws.Range("R10").Value = 47.039543820066
$ws.Range("S10").Value = 0.1761349325807084
$ws.Range("T10").Value = 0.1761349325807084

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Wnt7b"
$ws.Range("C11").Value = "Fzd4"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.420527
$ws.Range("H11").Value = 1.261581
$ws.Range("I11").Value = 0.852224964839111
$ws.Range("J11").Value = 0.852224964839111
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.251329
$ws.Range("N11").Value = 0.753987
$ws.Range("O11").Value = 0.004179334592343558
$ws.Range("P11").Value = 0.004179334592343557
$ws.Range("Q11").Value = 0.105690630383
$ws.Range("R11").Value = 0.951215673447
$ws.Range("S11").Value = 0.003561733276010869
$ws.Range("T11").Value = 0.003561733276010868
